# Update Multiple Language Button
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 5-9 (column A = Code, column B = English text, column C = Vietnamese text)
# Row 5: LOGIN
$ws.Range("B5").Value = "BG-Login.png"
$ws.Range("C5").Value = "BG-Login-VN.png"

# Row 6: CREATE
$ws.Range("B6").Value = "BG-Create_account.png"
$ws.Range("C6").Value = "BG-Create_account-VN.png"

# Row 7: BG-TITLE
$ws.Range("B7").Value = "BG-Title.png"
$ws.Range("C7").Value = "BG-Title-VN.png"

# Row 8: BG-PLAY -> BTN-PLAY
$ws.Range("A8").Value = "BTN-PLAY"
$ws.Range("B8").Value = "BTN-Play.png"
$ws.Range("C8").Value = "BTN-Play-VN.png"

# Row 9: BG-MINIGAME -> BTN-MINIGAME
$ws.Range("A9").Value = "BTN-MINIGAME"
$ws.Range("B9").Value = "BTN-Minigame.png"
$ws.Range("C9").Value = "BTN-Minigame.png"

# Update the selected cell to match the saved selection in the sheet view
$ws.Range("C15").Select()

$wb.Save()
